$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New column F header: BFO_DEF, styled like the other header cells (B1:E1)
$ws.Range("F1").Value2 = "BFO_DEF"
$ws.Range("F1").Font.Bold = $true
$ws.Range("F1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("F1").VerticalAlignment = -4160     # xlTop
$ws.Range("F1").Borders.LineStyle = 1         # xlContinuous (thin box border)

# BFO definition text ('BFO_DEF') for each mapping row, in the same order as the rows
# (rows 2-36 correspond 1:1 with BFO classes BFO_0000040 .. BFO_0000147)
$definitions = @(
    '[''A material entity is an independent continuant that has some portion of matter as proper or improper continuant part. [BFO]'']',
    '[''p is a process if p is an occurrent that has temporal proper parts and for some time t, p specifically depends on some material entity at t. [BFO]'', locstr("Process, i.e., a physical entity with a temporal evolution that ''has a meaning for the ontologist''", ''en'')]',
    '[''p is a process boundary if p is a temporal part of a process and p has no proper temporal parts. [BFO]'']',
    '[''A function is a disposition that exists in virtue of the bearer’s physical make-up and this physical make-up is something the bearer possesses because it came into being, either through evolution (in the case of natural biological entities) or through intentional design (in the case of artifacts), in order to realize processes of a certain sort. [BFO]'']',
    '[''A quality is a specifically dependent continuant that, in contrast to roles and dispositions, does not require any further process in order to be realized. [BFO]'']',
    '[''B is a process profile if there is some process c such that b is process profile of c. B is process profile of c holds when b is a proper occurrent part of c and there is some proper occurrent part d of c which has no parts in common with b and is mutually dependent on b and is such that b, c and d occupy the same temporal region. [BFO]'']',
    '[''An occurrent is an entity that unfolds itself in time or it is the instantaneous boundary of such an entity (for example a beginning or an ending) or it is a temporal or spatiotemporal region which such an entity occupies a temporal region or occupies a spatiotemporal region. [BFO]'']',
    '[''To say that b is a realizable entity is to say that b is a specifically dependent continuant that inheres in some independent continuant which is not a spatial region and is of a type instances of which are realized in processes of a correlated type.´[BFO]'', ''To say that b is a realizable entity is to say that b is a specifically dependent continuant that inheres in some independent continuant which is not a spatial region and is of a type instances of which are realized in processes of a correlated type. (axiom label in BFO2 Reference: [058-002])'']',
    '[''A continuant is an entity that persists, endures, or continues to exist through time while maintaining its identity. [BFO]'']',
    '[''b is an independent continuant if b is a continuant which is such that there is no c and no t such that b s-depends_on c at t. [BFO]'']',
    '[''A spatial region is a continuant entity that is a continuant part of space R as defined relative to some frame R. [BFO]'']',
    '[''An entity is anything that exists or has existed or will exist. [BFO]'']',
    '[''A temporal region is an occurrent entity that is part of time as defined relative to some reference frame. [BFO]'']',
    '[''A spatiotemporal region is an occurrent entity that is part of spacetime. [BFO]'']',
    '[''b is a specifically dependent continuant if b is a continuant and there is some independent continuant c which is not a spatial region and which is such that b specifically depends on c at every time t during the course of b’s existence. [BFO]'']',
    '[''b is a generically dependent continuant if b is a continuant that generically depends on one or more other entities. [BFO]'']',
    '[''B is a role means: b is a realizable entity and b exists because there is some single bearer that is in some special physical, social, or institutional set of circumstances in which this bearer does not have to be and b is not such that, if it ceases to exist, then the physical make-up of the bearer is thereby changed. [BFO]'']',
    '[''B is a disposition means: b is a realizable entity and b’s bearer is some material entity and b is such that if it ceases to exist, then its bearer is physically changed, and b’s realization occurs when and because this bearer is in some special physical circumstances, and this realization occurs in virtue of the bearer’s physical make-up. [BFO]'']',
    '[''B is a fiat object part if b is a material entity which is such that for all times t, if b exists at t then there is some object c such that b is a proper continuant part of c at t and c is demarcated from the remainder of c by a two-dimensional continuant fiat boundary. [BFO]'']',
    '[''B is a site means: b is a three-dimensional immaterial entity that is (partially or wholly) bounded by a material entity or it is a three-dimensional immaterial part thereof. [BFO]'']',
    '[''B is an object means: b is a material entity which manifests causal unity of one or other of the types causal unities and is of a type (a material universal) instances of which are maximal relative to this criterion of causal unity.\n\nTo say that b is causally unified means: b is a material entity which is such that its material parts are tied together in such a way that, in environments typical for entities of the type in question,if c, a continuant part of b that is in the interior of b at t, is larger than a certain threshold size (which will be determined differently from case to case, depending on factors such as porosity of external cover) and is moved in space to be at t at a location on the exterior of the spatial region that had been occupied by b at t, then either b’s other parts will be moved in coordinated fashion or b will be damaged (be affected, for example, by breakage or tearing) in the interval between t and t. [BFO]'']',
    '[''B is an object aggregate means: b is a material entity consisting exactly of a plurality of objects as member parts at all times at which b exists. [BFO]'']',
    '[''B is a relational quality if for some independent continuants c, d and for some time t: b is quality of c at t and b is quality of d at t. [BFO]'']',
    '[''A history is a process that is the sum of the totality of processes taking place in the spatiotemporal region occupied by a material entity or site, including processes on the surface of the entity or within the cavities to which it serves as host. [BFO]'']',
    '[''A zero-dimensional temporal region is a temporal region that is without extent. [BFO]'']',
    '[''A one-dimensional temporal region is a temporal region that is extended. [BFO]'']',
    '[]',
    '[''B is a continuant fiat boundary if b is an immaterial entity that is of zero, one or two dimensions and does not include a spatial region as part. [BFO]'']',
    '[''A two-dimensional spatial region is a spatial region that is of two dimensions. [BFO]'']',
    '[''A three-dimensional spatial region is a spatial region that is of three dimensions. [BFO]'']',
    '[''A zero-dimensional spatial region is a point in space. [BFO]'']',
    '[''A one-dimensional spatial region is a line or aggregate of lines stretching from one point in space to another. [BFO]'']',
    '[''A one-dimensional continuant fiat boundary is a continuous fiat line whose location is defined in relation to some material entity.  [BFO]'']',
    '[''A two-dimensional continuant fiat boundary (surface) is a self-connected fiat surface whose location is defined in relation to some material entity. [BFO]'']',
    '[''A zero-dimensional continuant fiat boundary is a fiat point whose location is defined in relation to some material entity. [BFO]'']'
)

for ($i = 0; $i -lt $definitions.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value2 = $definitions[$i]
}

# Row 37 (AFSP_0000001 -> BFO_0000029 'site') reuses the BFO_0000029 definition (same as row 21)
$ws.Cells.Item(37, 6).Value2 = $ws.Cells.Item(21, 6).Value2
